$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert two new data rows at the top of the data block
# (row 39), pushing the existing rows (39..144) down to (41..146).
$ws.Rows.Item(39).Resize(2).EntireRow.Insert()

# New row 39 ("Primera")
$ws.Cells.Item(39, 1).Value = 9
$ws.Cells.Item(39, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(39, 3).Value = "Metropolitana"
$ws.Cells.Item(39, 4).Value = 44497
$ws.Cells.Item(39, 5).Value = 13
$ws.Cells.Item(39, 6).Value = 100112017
$ws.Cells.Item(39, 7).Value = "Apio"
$ws.Cells.Item(39, 8).Value = "Americana (o)"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 50
$ws.Cells.Item(39, 11).Value = 7000
$ws.Cells.Item(39, 12).Value = 7000
$ws.Cells.Item(39, 13).Value = 7000
$ws.Cells.Item(39, 14).Value = "`$/docena de matas"
$ws.Cells.Item(39, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(39, 16).Value = 1167
$ws.Cells.Item(39, 17).Value = 6
$ws.Cells.Item(39, 18).Value = "Hortaliza"

# New row 40 ("Segunda")
$ws.Cells.Item(40, 1).Value = 9
$ws.Cells.Item(40, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(40, 3).Value = "Metropolitana"
$ws.Cells.Item(40, 4).Value = 44497
$ws.Cells.Item(40, 5).Value = 13
$ws.Cells.Item(40, 6).Value = 100112017
$ws.Cells.Item(40, 7).Value = "Apio"
$ws.Cells.Item(40, 8).Value = "Americana (o)"
$ws.Cells.Item(40, 9).Value = "Segunda"
$ws.Cells.Item(40, 10).Value = 80
$ws.Cells.Item(40, 11).Value = 5000
$ws.Cells.Item(40, 12).Value = 6000
$ws.Cells.Item(40, 13).Value = 5438
$ws.Cells.Item(40, 14).Value = "`$/docena de matas"
$ws.Cells.Item(40, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(40, 16).Value = 906
$ws.Cells.Item(40, 17).Value = 6
$ws.Cells.Item(40, 18).Value = "Hortaliza"
